$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 5 (A5:H5) down into the new row 6 (A6:H6)
$ws.Range("A5:H5").Copy()
$ws.Range("A6:H6").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = 5
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = "yes"
$ws.Range("G6").Value = 101
$ws.Range("H6").Value = "comment5"

$ws.Range("C9").Select()
